$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A12").Value = "LFBCJ8"
$ws.Range("B12").Value = "Buje de rodillo superior para Kyocera"
$ws.Range("C12").Value = "FS 1300 1110 1100 1024 1124 1128 1028 M2810 KM2820 ECOSYS M2030DN M2530 M2035 M2535 P2035 P2135"
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 80000
$ws.Range("F12").Value = 5
$ws.Range("G12").Value = 0
$ws.Range("H12").Formula = "=(E12-D12)*G12"
$ws.Range("I12").Formula = "=D12*F12"
$ws.Range("J12").Value = 0
